# Auto-generated edit script applying numeric updates per the commit diff.
# Updates profit-calculation columns (H-N) across multiple item rows on several sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 211.47058
$ws.Range("I33").Value = 199.24
$ws.Range("J33").Value = 245.44444
$ws.Range("K33").Value = 199.24
$ws.Range("L33").Value = 245.44444
$ws.Range("M33").Value = 29.75999999999999
$ws.Range("N33").Value = -703.44444
$ws.Range("H53").Value = 317.5
$ws.Range("I53").Value = 71.40000000000001
$ws.Range("J53").Value = 412.15384
$ws.Range("K53").Value = 71.40000000000001
$ws.Range("L53").Value = 412.15384
$ws.Range("M53").Value = 565.6
$ws.Range("N53").Value = -1686.15384
$ws.Range("H87").Value = 19821.54
$ws.Range("J87").Value = 19821.54
$ws.Range("L87").Value = 19821.54
$ws.Range("N87").Value = -22317.54
$ws.Range("H90").Value = 19821.54
$ws.Range("J90").Value = 19821.54
$ws.Range("L90").Value = 59464.62
$ws.Range("N90").Value = -71944.62
$ws.Range("H113").Value = 2643.8235
$ws.Range("I113").Value = 2556.3333
$ws.Range("J113").Value = 3300
$ws.Range("K113").Value = 2556.3333
$ws.Range("L113").Value = 3300
$ws.Range("M113").Value = 697.6667000000002
$ws.Range("N113").Value = -9808
$ws.Range("H127").Value = 1665.7273
$ws.Range("I127").Value = 637.5
$ws.Range("J127").Value = 2899.6
$ws.Range("K127").Value = 1912.5
$ws.Range("L127").Value = 8698.799999999999
$ws.Range("M127").Value = 3047.5
$ws.Range("N127").Value = -18618.8
$ws.Range("H132").Value = 1588.025
$ws.Range("I132").Value = 1592.5
$ws.Range("J132").Value = 1503
$ws.Range("K132").Value = 4777.5
$ws.Range("L132").Value = 4509
$ws.Range("M132").Value = -2247.5
$ws.Range("N132").Value = -9569

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3478.08
$ws.Range("I32").Value = 2505.5173
$ws.Range("J32").Value = 9986.77
$ws.Range("K32").Value = 2505.5173
$ws.Range("L32").Value = 9986.77
$ws.Range("M32").Value = -2218.5173
$ws.Range("N32").Value = -10560.77
$ws.Range("H45").Value = 1474.2
$ws.Range("J45").Value = 1431.75
$ws.Range("L45").Value = 1431.75
$ws.Range("N45").Value = -2185.75
$ws.Range("H53").Value = 15470.286
$ws.Range("I53").Value = 3659.75
$ws.Range("J53").Value = 31217.666
$ws.Range("K53").Value = 3659.75
$ws.Range("L53").Value = 31217.666
$ws.Range("M53").Value = -2977.75
$ws.Range("N53").Value = -32581.666
$ws.Range("H122").Value = 6312
$ws.Range("I122").Value = 9472.154
$ws.Range("J122").Value = 2577.2727
$ws.Range("K122").Value = 28416.462
$ws.Range("L122").Value = 7731.8181
$ws.Range("M122").Value = -25966.462
$ws.Range("N122").Value = -12631.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 40774
$ws.Range("J51").Value = 40774
$ws.Range("L51").Value = 40774
$ws.Range("N51").Value = -41756

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27898.426
$ws.Range("I31").Value = 38522.785
$ws.Range("J31").Value = 3108.25
$ws.Range("K31").Value = 38522.785
$ws.Range("L31").Value = 3108.25
$ws.Range("M31").Value = -38227.785
$ws.Range("N31").Value = -3698.25
$ws.Range("H34").Value = 27898.426
$ws.Range("I34").Value = 38522.785
$ws.Range("J34").Value = 3108.25
$ws.Range("K34").Value = 38522.785
$ws.Range("L34").Value = 3108.25
$ws.Range("M34").Value = -38320.785
$ws.Range("N34").Value = -3512.25
$ws.Range("H58").Value = 1028.0625
$ws.Range("I58").Value = 937.0417
$ws.Range("J58").Value = 1301.125
$ws.Range("K58").Value = 937.0417
$ws.Range("L58").Value = 1301.125
$ws.Range("M58").Value = -734.0417
$ws.Range("N58").Value = -1707.125
$ws.Range("H94").Value = 167636.3
$ws.Range("I94").Value = 200543.2
$ws.Range("J94").Value = 144131.36
$ws.Range("K94").Value = 200543.2
$ws.Range("L94").Value = 144131.36
$ws.Range("M94").Value = -200092.2
$ws.Range("N94").Value = -145033.36
$ws.Range("H136").Value = 1028.0625
$ws.Range("I136").Value = 937.0417
$ws.Range("J136").Value = 1301.125
$ws.Range("K136").Value = 2811.1251
$ws.Range("L136").Value = 3903.375
$ws.Range("M136").Value = -261.1251000000002
$ws.Range("N136").Value = -9003.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 369.85715
$ws.Range("I107").Value = 216.5
$ws.Range("J107").Value = 431.2
$ws.Range("K107").Value = 649.5
$ws.Range("L107").Value = 1293.6
$ws.Range("M107").Value = 1270.5
$ws.Range("N107").Value = -5133.6
$ws.Range("H113").Value = 827.675
$ws.Range("I113").Value = 1307.7646
$ws.Range("J113").Value = 472.82608
$ws.Range("K113").Value = 3923.2938
$ws.Range("L113").Value = 1418.47824
$ws.Range("M113").Value = -1753.2938
$ws.Range("N113").Value = -5758.47824

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 42499.5
$ws.Range("J47").Value = 42499.5
$ws.Range("L47").Value = 42499.5
$ws.Range("N47").Value = -43479.5
$ws.Range("H52").Value = 42499.5
$ws.Range("J52").Value = 42499.5
$ws.Range("L52").Value = 42499.5
$ws.Range("N52").Value = -42965.5
$ws.Range("H55").Value = 438.13333
$ws.Range("I55").Value = 407.77777
$ws.Range("J55").Value = 483.66666
$ws.Range("K55").Value = 407.77777
$ws.Range("L55").Value = 483.66666
$ws.Range("M55").Value = -234.77777
$ws.Range("N55").Value = -829.66666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 25200
$ws.Range("I29").Value = 29000
$ws.Range("J29").Value = 22666.666
$ws.Range("K29").Value = 29000
$ws.Range("L29").Value = 22666.666
$ws.Range("M29").Value = -28710
$ws.Range("N29").Value = -23246.666
$ws.Range("H43").Value = 55000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 55000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 55000
$ws.Range("M43").Value = ""
$ws.Range("N43").Value = -55298
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").Value = ""
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""
$ws.Range("H126").Value = 589260.5600000001
$ws.Range("I126").Value = 770047.7
$ws.Range("K126").Value = 2310143.1
$ws.Range("M126").Value = -2307673.1
$ws.Range("H132").Value = 579.24194
$ws.Range("I132").Value = 395.40817
$ws.Range("J132").Value = 1272.1538
$ws.Range("K132").Value = 1186.22451
$ws.Range("L132").Value = 3816.4614
$ws.Range("M132").Value = 1343.77549
$ws.Range("N132").Value = -8876.4614
$ws.Range("H136").Value = 354.95123
$ws.Range("I136").Value = 290.35135
$ws.Range("J136").Value = 952.5
$ws.Range("K136").Value = 871.0540500000001
$ws.Range("L136").Value = 2857.5
$ws.Range("M136").Value = 1678.94595
$ws.Range("N136").Value = -7957.5

